# scrum.xlsx update: refresh Product Backlog status, extend several
# Sprint Backlog rows with "Effort Actual" / "Status" data, add a new
# Sprint 2 task row for Vaadin navigation, and refresh view selections.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Product Backlog")
$ws2 = $wb.Worksheets.Item("Sprint Backlog")

# ---------------------------------------------------------------------
# Product Backlog: "Conduct a patient session" is now in progress
# ---------------------------------------------------------------------
$ws1.Range("H3").Value = "Work In Progress"

# ---------------------------------------------------------------------
# Sprint Backlog: fill in effort-actual / status columns for several
# sprint-2 rows (new shared strings must be introduced in this exact
# order: "10h" then "Figure out navigation for Vaadin Framework" then
# "Basic GUI, navigation between pages" then "?").
# ---------------------------------------------------------------------

# Row 12 - Patient session information
$ws2.Range("I12").Value = "10h"
$ws2.Range("J12").Value = "8h"

# Row 14 - Test Data for development
$ws2.Range("J14").Value = "2h"
$ws2.Range("K14").Value = "Done"

# Row 15 - MedicationService
$ws2.Range("I15").Value = "2h"
$ws2.Range("J15").Value = "4h"
$ws2.Range("K15").Value = "Work In Progress"

# Row 16 - Medication prescription
$ws2.Range("J16").Value = "0h"
$ws2.Range("K16").Value = "Work In Progress"

# Row 17 - Code documentation (also story-point estimate correction)
$ws2.Range("A17").Value = 0.1
$ws2.Range("J17").Value = "0h"
$ws2.Range("K17").Value = "Work In Progress"

# Row 18 - Test data prep for task 2.2
$ws2.Range("J18").Value = "0.5h"
$ws2.Range("K18").Value = "Done"

# Row 19 (new) - Vaadin Framework navigation task
$ws2.Range("A19").Value = 0.2
$ws2.Range("B19").Value = 2
$ws2.Range("D19").Value = "Figure out navigation for Vaadin Framework"
$ws2.Range("C19").Value = "Basic GUI, navigation between pages"
$ws2.Range("E19").Value = "UI"
$ws2.Range("F19").Value = "Meyer"
$ws2.Range("G19").Value = "low"
$ws2.Range("H19").Value = "?"
$ws2.Range("I19").Value = "3h"
$ws2.Range("J19").Value = "1h"
$ws2.Range("K19").Value = "Work In Progress"

# ---------------------------------------------------------------------
# View state: refresh selections on both sheets
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H5").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("J16").Select()

# Book window geometry
$excel.Windows.Item(1).Width = 28800
$excel.Windows.Item(1).Height = 17560
